$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 310
$ws1.Range("F3").Value = 508
$ws1.Range("F5").Value = 8750
$ws1.Range("F7").Value = 11141
$ws1.Range("F15").Value = 301
$ws1.Range("F22").Value = 1899
$ws1.Range("F23").Value = 708
$ws1.Range("F26").Value = 293
$ws1.Range("F28").Value = 606
$ws1.Range("F30").Value = 1295
$ws1.Range("F31").Value = 27
$ws1.Range("F38").Value = 352
$ws1.Range("F39").Value = 303
$ws1.Range("F41").Value = 144
$ws1.Range("F42").Value = 535
$ws1.Range("F43").Value = 384
$ws1.Range("F46").Value = 657
$ws1.Range("F48").Value = 155

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 31
$ws2.Range("G19").Value = "不可售"

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 508
$ws4.Range("F8").Value = 8750
$ws4.Range("F10").Value = 11141
$ws4.Range("F15").Value = 301
$ws4.Range("F19").Value = 1899
$ws4.Range("F20").Value = 708
$ws4.Range("F23").Value = 293
$ws4.Range("F25").Value = 606
$ws4.Range("F29").Value = 1295
$ws4.Range("F30").Value = 27
$ws4.Range("F33").Value = 31
$ws4.Range("F40").Value = 352
$ws4.Range("F41").Value = 535
$ws4.Range("F42").Value = 384
$ws4.Range("F46").Value = 657
$ws4.Range("F48").Value = 155
